$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5024.9165  # H40
$ws.Cells.Item(40, 9).Value = 3300  # I40
$ws.Cells.Item(40, 10).Value = 5887.375  # J40
$ws.Cells.Item(40, 11).Value = 3300  # K40
$ws.Cells.Item(40, 12).Value = 5887.375  # L40
$ws.Cells.Item(40, 13).Value = -3125  # M40
$ws.Cells.Item(40, 14).Value = -6237.375  # N40
$ws.Cells.Item(41, 8).Value = 3193.389  # H41
$ws.Cells.Item(41, 10).Value = 3065.2856  # J41
$ws.Cells.Item(41, 12).Value = 3065.2856  # L41
$ws.Cells.Item(41, 14).Value = -3945.2856  # N41
$ws.Cells.Item(43, 8).Value = 232817.56  # H43
$ws.Cells.Item(43, 9).Value = 7916.6665  # I43
$ws.Cells.Item(43, 10).Value = 345268  # J43
$ws.Cells.Item(43, 11).Value = 7916.6665  # K43
$ws.Cells.Item(43, 12).Value = 345268  # L43
$ws.Cells.Item(43, 13).Value = -7847.6665  # M43
$ws.Cells.Item(43, 14).Value = -345406  # N43
$ws.Cells.Item(62, 8).Value = 45500  # H62
$ws.Cells.Item(62, 9).Value = 0  # I62
$ws.Cells.Item(62, 11).Value = 0  # K62
$ws.Cells.Item(62, 13).ClearContents()  # M62
$ws.Cells.Item(65, 8).Value = 45500  # H65
$ws.Cells.Item(65, 9).Value = 0  # I65
$ws.Cells.Item(65, 11).Value = 0  # K65
$ws.Cells.Item(65, 13).ClearContents()  # M65
$ws.Cells.Item(98, 8).Value = 62500524  # H98
$ws.Cells.Item(98, 9).Value = 62500524  # I98
$ws.Cells.Item(98, 11).Value = 62500524  # K98
$ws.Cells.Item(98, 13).Value = -62499026  # M98
$ws.Cells.Item(113, 8).Value = 45458180  # H113
$ws.Cells.Item(113, 9).Value = 14290185  # I113
$ws.Cells.Item(113, 10).Value = 100002180  # J113
$ws.Cells.Item(113, 11).Value = 14290185  # K113
$ws.Cells.Item(113, 12).Value = 100002180  # L113
$ws.Cells.Item(113, 13).Value = -14286931  # M113
$ws.Cells.Item(113, 14).Value = -100008688  # N113
$ws.Cells.Item(122, 8).Value = 62500524  # H122
$ws.Cells.Item(122, 9).Value = 62500524  # I122
$ws.Cells.Item(122, 11).Value = 187501572  # K122
$ws.Cells.Item(122, 13).Value = -187499122  # M122
$ws.Cells.Item(126, 8).Value = 74110  # H126
$ws.Cells.Item(126, 10).Value = 72775  # J126
$ws.Cells.Item(126, 12).Value = 72775  # L126
$ws.Cells.Item(126, 14).Value = -82655  # N126
$ws.Cells.Item(128, 8).Value = 64666.668  # H128
$ws.Cells.Item(128, 10).Value = 64666.668  # J128
$ws.Cells.Item(128, 12).Value = 64666.668  # L128
$ws.Cells.Item(128, 14).Value = -74626.66800000001  # N128
$ws.Cells.Item(135, 8).Value = 3546  # H135
$ws.Cells.Item(135, 9).Value = 3546  # I135
$ws.Cells.Item(135, 11).Value = 31914  # K135
$ws.Cells.Item(135, 13).Value = -29379  # M135
$ws.Cells.Item(138, 8).Value = 3235.4353  # H138
$ws.Cells.Item(138, 10).Value = 3357.3635  # J138
$ws.Cells.Item(138, 12).Value = 10072.0905  # L138
$ws.Cells.Item(138, 14).Value = -20352.0905  # N138
$ws.Cells.Item(141, 8).Value = 3216.5715  # H141
$ws.Cells.Item(141, 9).Value = 3160.8  # I141
$ws.Cells.Item(141, 11).Value = 9482.400000000001  # K141
$ws.Cells.Item(141, 13).Value = -4302.400000000001  # M141

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1571.875  # H4
$ws.Cells.Item(4, 9).Value = 1137.8334  # I4
$ws.Cells.Item(4, 11).Value = 1137.8334  # K4
$ws.Cells.Item(4, 13).Value = -1021.8334  # M4
$ws.Cells.Item(32, 8).Value = 13160965  # H32
$ws.Cells.Item(32, 9).Value = 15153353  # I32
$ws.Cells.Item(32, 10).Value = 11205.4  # J32
$ws.Cells.Item(32, 11).Value = 15153353  # K32
$ws.Cells.Item(32, 12).Value = 11205.4  # L32
$ws.Cells.Item(32, 13).Value = -15153066  # M32
$ws.Cells.Item(32, 14).Value = -11779.4  # N32
$ws.Cells.Item(88, 8).Value = 1681.3478  # H88
$ws.Cells.Item(88, 9).Value = 1293.2  # I88
$ws.Cells.Item(88, 10).Value = 1979.9231  # J88
$ws.Cells.Item(88, 11).Value = 1293.2  # K88
$ws.Cells.Item(88, 12).Value = 1979.9231  # L88
$ws.Cells.Item(88, 13).Value = -887.2  # M88
$ws.Cells.Item(88, 14).Value = -2791.9231  # N88
$ws.Cells.Item(91, 8).Value = 1681.3478  # H91
$ws.Cells.Item(91, 9).Value = 1293.2  # I91
$ws.Cells.Item(91, 10).Value = 1979.9231  # J91
$ws.Cells.Item(91, 11).Value = 1293.2  # K91
$ws.Cells.Item(91, 12).Value = 1979.9231  # L91
$ws.Cells.Item(91, 13).Value = 110.8  # M91
$ws.Cells.Item(91, 14).Value = -4787.9231  # N91
$ws.Cells.Item(102, 8).Value = 5622  # H102
$ws.Cells.Item(102, 9).Value = 5757.294  # I102
$ws.Cells.Item(102, 11).Value = 5757.294  # K102
$ws.Cells.Item(102, 13).Value = -4135.294  # M102
$ws.Cells.Item(110, 8).Value = 2204.1428  # H110
$ws.Cells.Item(110, 10).Value = 0  # J110
$ws.Cells.Item(110, 12).Value = 0  # L110
$ws.Cells.Item(110, 14).ClearContents()  # N110
$ws.Cells.Item(112, 8).Value = 71212.25  # H112
$ws.Cells.Item(112, 10).Value = 71212.25  # J112
$ws.Cells.Item(112, 12).Value = 71212.25  # L112
$ws.Cells.Item(112, 14).Value = -74166.25  # N112

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 84996.664  # H2
$ws.Cells.Item(2, 10).Value = 84996.664  # J2
$ws.Cells.Item(2, 12).Value = 84996.664  # L2
$ws.Cells.Item(2, 14).Value = -85222.664  # N2
$ws.Cells.Item(97, 8).Value = 24484.666  # H97
$ws.Cells.Item(97, 9).Value = 5278.6  # I97
$ws.Cells.Item(97, 10).Value = 48492.25  # J97
$ws.Cells.Item(97, 11).Value = 5278.6  # K97
$ws.Cells.Item(97, 12).Value = 48492.25  # L97
$ws.Cells.Item(97, 13).Value = -4287.6  # M97
$ws.Cells.Item(97, 14).Value = -50474.25  # N97
$ws.Cells.Item(105, 8).Value = 2702.8  # H105
$ws.Cells.Item(105, 9).Value = 2702.375  # I105
$ws.Cells.Item(105, 11).Value = 2702.375  # K105
$ws.Cells.Item(105, 13).Value = -955.375  # M105

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1225.9584  # H58
$ws.Cells.Item(58, 9).Value = 1119.9474  # I58
$ws.Cells.Item(58, 10).Value = 1628.8  # J58
$ws.Cells.Item(58, 11).Value = 1119.9474  # K58
$ws.Cells.Item(58, 12).Value = 1628.8  # L58
$ws.Cells.Item(58, 13).Value = -916.9474  # M58
$ws.Cells.Item(58, 14).Value = -2034.8  # N58
$ws.Cells.Item(105, 8).Value = 2032.2222  # H105
$ws.Cells.Item(105, 9).Value = 1974.6471  # I105
$ws.Cells.Item(105, 11).Value = 1974.6471  # K105
$ws.Cells.Item(105, 13).Value = -227.6470999999999  # M105
$ws.Cells.Item(129, 8).Value = 84996.664  # H129
$ws.Cells.Item(129, 10).Value = 84996.664  # J129
$ws.Cells.Item(129, 12).Value = 84996.664  # L129
$ws.Cells.Item(129, 14).Value = -94996.664  # N129
$ws.Cells.Item(136, 8).Value = 1225.9584  # H136
$ws.Cells.Item(136, 9).Value = 1119.9474  # I136
$ws.Cells.Item(136, 10).Value = 1628.8  # J136
$ws.Cells.Item(136, 11).Value = 3359.8422  # K136
$ws.Cells.Item(136, 12).Value = 4886.4  # L136
$ws.Cells.Item(136, 13).Value = -809.8422  # M136
$ws.Cells.Item(136, 14).Value = -9986.4  # N136

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 49892.5  # H37
$ws.Cells.Item(37, 10).Value = 49892.5  # J37
$ws.Cells.Item(37, 12).Value = 149677.5  # L37
$ws.Cells.Item(37, 14).Value = -149901.5  # N37
$ws.Cells.Item(55, 8).Value = 950  # H55
$ws.Cells.Item(55, 9).Value = 950  # I55
$ws.Cells.Item(55, 11).Value = 2850  # K55
$ws.Cells.Item(55, 13).Value = -2673  # M55
$ws.Cells.Item(137, 8).Value = 4682.6  # H137
$ws.Cells.Item(137, 10).Value = 6610.2144  # J137
$ws.Cells.Item(137, 12).Value = 19830.6432  # L137
$ws.Cells.Item(137, 14).Value = -30030.6432  # N137

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 13100.6  # H20
$ws.Cells.Item(20, 10).Value = 13166  # J20
$ws.Cells.Item(20, 12).Value = 13166  # L20
$ws.Cells.Item(20, 14).Value = -13656  # N20
$ws.Cells.Item(24, 8).Value = 4998  # H24
$ws.Cells.Item(24, 10).Value = 4998  # J24
$ws.Cells.Item(24, 12).Value = 4998  # L24
$ws.Cells.Item(24, 14).Value = -5344  # N24
$ws.Cells.Item(111, 8).Value = 81045.75  # H111
$ws.Cells.Item(111, 10).Value = 81045.75  # J111
$ws.Cells.Item(111, 12).Value = 81045.75  # L111
$ws.Cells.Item(111, 14).Value = -87179.75  # N111
$ws.Cells.Item(122, 8).Value = 1394.9474  # H122
$ws.Cells.Item(122, 9).Value = 1156.7858  # I122
$ws.Cells.Item(122, 11).Value = 3470.3574  # K122
$ws.Cells.Item(122, 13).Value = -1020.3574  # M122

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 51995  # H6
$ws.Cells.Item(6, 10).Value = 51995  # J6
$ws.Cells.Item(6, 12).Value = 51995  # L6
$ws.Cells.Item(6, 14).Value = -52219  # N6
$ws.Cells.Item(46, 8).Value = 4821.909  # H46
$ws.Cells.Item(46, 9).Value = 2181.7334  # I46
$ws.Cells.Item(46, 10).Value = 10479.429  # J46
$ws.Cells.Item(46, 11).Value = 2181.7334  # K46
$ws.Cells.Item(46, 12).Value = 10479.429  # L46
$ws.Cells.Item(46, 13).Value = -1993.7334  # M46
$ws.Cells.Item(46, 14).Value = -10855.429  # N46
$ws.Cells.Item(61, 8).Value = 2138.25  # H61
$ws.Cells.Item(61, 10).Value = 2994  # J61
$ws.Cells.Item(61, 12).Value = 2994  # L61
$ws.Cells.Item(61, 14).Value = -3398  # N61
$ws.Cells.Item(110, 8).Value = 45040.832  # H110
$ws.Cells.Item(110, 10).Value = 45040.832  # J110
$ws.Cells.Item(110, 12).Value = 45040.832  # L110
$ws.Cells.Item(110, 14).Value = -53220.832  # N110
$ws.Cells.Item(113, 8).Value = 2138.25  # H113
$ws.Cells.Item(113, 10).Value = 2994  # J113
$ws.Cells.Item(113, 12).Value = 2994  # L113
$ws.Cells.Item(113, 14).Value = -7334  # N113
$ws.Cells.Item(122, 8).Value = 6968.115  # H122
$ws.Cells.Item(122, 9).Value = 7201.8237  # I122
$ws.Cells.Item(122, 11).Value = 21605.4711  # K122
$ws.Cells.Item(122, 13).Value = -19155.4711  # M122
$ws.Cells.Item(132, 8).Value = 911763.9399999999  # H132
$ws.Cells.Item(132, 9).Value = 2710.7778  # I132
$ws.Cells.Item(132, 10).Value = 5002503  # J132
$ws.Cells.Item(132, 11).Value = 8132.3334  # K132
$ws.Cells.Item(132, 12).Value = 15007509  # L132
$ws.Cells.Item(132, 13).Value = -5602.3334  # M132
$ws.Cells.Item(132, 14).Value = -15012569  # N132

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 62790.363  # H4
$ws.Cells.Item(4, 10).Value = 30061.875  # J4
$ws.Cells.Item(4, 12).Value = 30061.875  # L4
$ws.Cells.Item(4, 14).Value = -30287.875  # N4
$ws.Cells.Item(50, 8).Value = 5000  # H50
$ws.Cells.Item(50, 9).Value = 4000  # I50
$ws.Cells.Item(50, 11).Value = 4000  # K50
$ws.Cells.Item(50, 13).Value = -3369  # M50
$ws.Cells.Item(61, 8).Value = 19833  # H61
$ws.Cells.Item(61, 9).Value = 9500  # I61
$ws.Cells.Item(61, 10).Value = 40499  # J61
$ws.Cells.Item(61, 11).Value = 9500  # K61
$ws.Cells.Item(61, 12).Value = 40499  # L61
$ws.Cells.Item(61, 13).Value = -9208  # M61
$ws.Cells.Item(61, 14).Value = -41083  # N61
$ws.Cells.Item(75, 8).Value = 4079993.5  # H75
$ws.Cells.Item(75, 10).Value = 99992  # J75
$ws.Cells.Item(75, 12).Value = 99992  # L75
$ws.Cells.Item(75, 14).Value = -101864  # N75
$ws.Cells.Item(78, 8).Value = 4079993.5  # H78
$ws.Cells.Item(78, 10).Value = 99992  # J78
$ws.Cells.Item(78, 12).Value = 299976  # L78
$ws.Cells.Item(78, 14).Value = -309336  # N78
$ws.Cells.Item(136, 8).Value = 2166.1667  # H136
$ws.Cells.Item(136, 9).Value = 1250.5  # I136
$ws.Cells.Item(136, 11).Value = 3751.5  # K136
$ws.Cells.Item(136, 13).Value = -1201.5  # M136
